$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comprados")
$ws.Activate()

# Remove the obsolete "Nac. - OPT 2" column (old column E). This shifts the
# old "Importados" column (F) left into E, updates the used dimension and
# row/col spans automatically.
$ws.Range("E1").EntireColumn.Delete()

# Update the approver assigned to several "Nacional" rows: SEBAREZE was
# replaced by either LUCIAE11 or ERIKBARB depending on the material.
$ws.Range("D3").Value = "LUCIAE11"
$ws.Range("D4").Value = "ERIKBARB"
$ws.Range("D5").Value = "LUCIAE11"
$ws.Range("D6").Value = "LUCIAE11"
$ws.Range("D7").Value = "ERIKBARB"
$ws.Range("D8").Value = "LUCIAE11"
$ws.Range("D13").Value = "LUCIAE11"
$ws.Range("D16").Value = "ERIKBARB"
$ws.Range("D17").Value = "ERIKBARB"
$ws.Range("D18").Value = "ERIKBARB"
$ws.Range("D19").Value = "ERIKBARB"
$ws.Range("D20").Value = "ERIKBARB"
$ws.Range("D21").Value = "ERIKBARB"
$ws.Range("D22").Value = "ERIKBARB"
$ws.Range("D23").Value = "ERIKBARB"
$ws.Range("D24").Value = "ERIKBARB"
$ws.Range("D25").Value = "ERIKBARB"

# Update the "last updated" notes to reflect the new alignment date.
$ws.Range("C1").Value = "Atualizado em: 22/11/2023"
$ws.Range("A30").Value = "Atualizado conforme alinhamento com o Carlos Sousa no dia 22/11/2023."

# Restore the selected cell shown when the workbook is reopened.
$ws.Range("D5").Select()
